$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 107.3665618896484
$ws.Range("C2").Value = 64

$ws.Range("B3").Value = 180.0963878631592
$ws.Range("C3").Value = 0

$ws.Range("B4").Value = 279.9334526062012

$ws.Range("B5").Value = 257.5397491455078

$ws.Range("B6").Value = 164.7019386291504

$ws.Range("B7").Value = 275.2518653869629

$ws.Range("B8").Value = 275.3086090087891

$ws.Range("B9").Value = 276.2391567230225

$ws.Range("B10").Value = 145.0092792510986

$ws.Range("B11").Value = 198.6324787139893

$ws.Range("B12").Value = 216.0079479217529

$ws.Range("B13").Value = 228.0861139297485
